$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column (H1), copying the formatting from the
# neighboring header cell (G1) so it gets the same style (bold, bordered,
# centered) as the rest of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the Save values for each data row
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
